$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new pair names (bestFit-style width ~17.75)
$ws.Columns.Item(1).ColumnWidth = 17

# Highlight the two new "btm_eth" / "eos_eth" pair names in column A (yellow fill),
# matching the highlight already applied to A8 ("etc_eth")
$ws.Range("A8").Interior.Color = 65535

# Row 13: btm_eth
$ws.Range("A13").Value = "btm_eth"
$ws.Range("A13").Interior.Color = 65535
$ws.Range("B13").Value = 0.00000001
$ws.Range("C13").Value = "ETH"
$ws.Range("D13").Value = 0.00000001
$ws.Range("E13").Value = "BTM"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = "BTM"
$ws.Range("H13").Value = "暂时没有"
$ws.Range("I13").Value = "暂时没有"
$ws.Range("C13").HorizontalAlignment = -4152
$ws.Range("E13").HorizontalAlignment = -4152
$ws.Range("G13").HorizontalAlignment = -4152
$ws.Range("H13").HorizontalAlignment = -4152
$ws.Range("I13").HorizontalAlignment = -4152

# Row 14: eos_eth
$ws.Range("A14").Value = "eos_eth"
$ws.Range("A14").Interior.Color = 65535
$ws.Range("B14").Value = 0.00000001
$ws.Range("C14").Value = "ETH"
$ws.Range("D14").Value = 0.00000001
$ws.Range("E14").Value = "EOS"
$ws.Range("F14").Value = 0.1
$ws.Range("G14").Value = "EOS"
$ws.Range("H14").Value = "暂时没有"
$ws.Range("I14").Value = "暂时没有"
$ws.Range("C14").HorizontalAlignment = -4152
$ws.Range("E14").HorizontalAlignment = -4152
$ws.Range("G14").HorizontalAlignment = -4152
$ws.Range("H14").HorizontalAlignment = -4152
$ws.Range("I14").HorizontalAlignment = -4152

# Row 15: eos_usdt
$ws.Range("A15").Value = "eos_usdt"
$ws.Range("B15").Value = 0.0001
$ws.Range("C15").Value = "USDT"
$ws.Range("D15").Value = 0.0001
$ws.Range("E15").Value = "EOS"
$ws.Range("F15").Value = 0.1
$ws.Range("G15").Value = "EOS"
$ws.Range("H15").Value = "暂时没有"
$ws.Range("I15").Value = "暂时没有"
$ws.Range("C15").HorizontalAlignment = -4152
$ws.Range("E15").HorizontalAlignment = -4152
$ws.Range("G15").HorizontalAlignment = -4152
$ws.Range("H15").HorizontalAlignment = -4152
$ws.Range("I15").HorizontalAlignment = -4152

# Row 17: note that eos_usdt needs crossing logic
$ws.Range("A17").Value = "需要Crossing Logic"
$ws.Range("A17").Interior.Color = 65535

$ws.Range("A17").Select()
